$wb = $excel.ActiveWorkbook

$renames = @{
    "m9-EX_cit_e-7.44-iModulon"     = "m9-EX_cit_e7.44-iModulon"
    "m9-EX_cit_e-7.44-Subsystem"    = "m9-EX_cit_e7.44-Subsystem"
    "m9-EX_fer_e-2.91-iModulon"     = "m9-EX_fer_e2.91-iModulon"
    "m9-EX_fer_e-2.91-Subsystem"    = "m9-EX_fer_e2.91-Subsystem"
    "m9-EX_glc_e-7.44-iModulon"     = "m9-EX_glc_e7.44-iModulon"
    "m9-EX_glc_e-7.44-Subsystem"    = "m9-EX_glc_e7.44-Subsystem"
    "m9-EX_ser__L_e-14.88-iModulon" = "m9-EX_ser__L_e14.88-iModulon"
    "m9-EX_ser__L_e-14.88-Subsystem"= "m9-EX_ser__L_e14.88-Subsystem"
}

foreach ($ws in $wb.Worksheets) {
    $oldName = $ws.Name
    if ($renames.ContainsKey($oldName)) {
        $ws.Name = $renames[$oldName]
    }
}
